# Update "File Extensions" reference sheet with additional CAD-tool columns
# (EAGLE CAM Processor 8.2.2 / CAMtastic7 / eagle version=7.6.0) and new rows
# for Top Paste, Bot Paste, Internal Signal Layers, Gerber Drawings and
# EIA NC Drill File, plus a few corrected/expanded extension lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "E"/"F"/"G" header columns on row 3
$ws.Range("E3").Value = "EAGLE CAM Processor 8.2.2"
$ws.Range("F3").Value = "CAMtastic7"
$ws.Range("G3").Value = "eagle version=7.6.0"

# Top Silkscreen (row 4)
$ws.Range("E4").Value = ".gts"
$ws.Range("F4").Value = ".gto"
$ws.Range("G4").Value = ".plc"

# Top Soldermask (row 5)
$ws.Range("E5").Value = ".gtm"
$ws.Range("F5").Value = ".gts"
$ws.Range("G5").Value = ".stc"

# Top Copper (row 6)
$ws.Range("E6").Value = ".gtl"
$ws.Range("F6").Value = ".gtl"
$ws.Range("G6").Value = ".cmp"

# Bottom Copper (row 7)
$ws.Range("E7").Value = ".gbl"
$ws.Range("F7").Value = ".gbl"
$ws.Range("G7").Value = ".sol"

# Bottom Soldermask (row 8)
$ws.Range("E8").Value = ".gbm"
$ws.Range("F8").Value = ".gbs"
$ws.Range("G8").Value = ".sts"

# Bottom Silkscreen (row 9)
$ws.Range("E9").Value = ".gbs"
$ws.Range("G9").Value = ".pls"

# NC Drill File (row 10)
$ws.Range("D10").Value = "thruhole.tap, .tap, .npt"
$ws.Range("E10").Value = ".drd, .txt"
$ws.Range("G10").Value = ".drd"

# Board Outline (row 11)
$ws.Range("D11").Value = ".bol"
$ws.Range("E11").Value = ".gbo"
$ws.Range("G11").Value = ".dim"

# Internal Layers -> Internal Plane Layers (row 12)
$ws.Range("A12").Value = "Internal Plane Layers"
$ws.Range("D12").Value = "in#"
$ws.Range("G12").Value = ".i#"

# Top Paste (row 15, new)
$ws.Range("A15").Value = "Top Paste"
$ws.Range("B15").Value = ".gtp"
$ws.Range("D15").Value = ".spt"
$ws.Range("E15").Value = ".gtp"
$ws.Range("F15").Value = ".gtp"
$ws.Range("G15").Value = ".crc"

# Bot Paste (row 16, new)
$ws.Range("A16").Value = "Bot Paste"
$ws.Range("B16").Value = ".gbp"
$ws.Range("D16").Value = ".spb"
$ws.Range("E16").Value = ".gbp"
$ws.Range("G16").Value = ".crs"

# Internal Signal Layers (row 17, new)
$ws.Range("A17").Value = "Internal Signal Layers"
$ws.Range("B17").Value = ".g1"
$ws.Range("D17").Value = ".in#"

# Gerber Drawings (row 18, new)
$ws.Range("A18").Value = "Gerber Drawings"
$ws.Range("B18").Value = ".gm#, .gd#, .gg#"
$ws.Range("D18").Value = ".drd, .asb, .ast, .fab"
$ws.Range("E18").Value = ".gtd"
$ws.Range("F18").Value = ".gm#"

# EIA NC Drill File (row 19, new)
$ws.Range("A19").Value = "EIA NC Drill File"
$ws.Range("B19").Value = ".drl"
